# aggregated_data_calculator: calculates average adjusted returns and average
# historical data. As part of this change the "Equity" and "Bond" sheets are
# swapped (their tab names and their row data are exchanged) and the Expense
# Ratio / Dividend Yield columns (F, G) are normalized from "percent-like"
# numbers (e.g. 1.62 meaning 1.62%) to true decimal fractions (e.g. 0.0162).

$wb = $excel.ActiveWorkbook

$wsEquitySlot = $wb.Worksheets.Item(1)   # physical sheet1.xml, currently named "Equity"
$wsAlternative = $wb.Worksheets.Item(2)  # physical sheet2.xml, "Alternative" (name unchanged)
$wsBondSlot = $wb.Worksheets.Item(3)     # physical sheet3.xml, currently named "Bond"

# ---------------------------------------------------------------------------
# 1) Swap the two sheet tab names (Equity <-> Bond). Use a temporary name to
#    avoid a duplicate-name collision while the swap is in progress.
# ---------------------------------------------------------------------------
$wsBondSlot.Name = "Equity_tmp"
$wsEquitySlot.Name = "Bond"
$wsBondSlot.Name = "Equity"

# From here on:
#   $wsEquitySlot (physical sheet1.xml) is now named "Bond" and will hold the bond fund rows.
#   $wsBondSlot   (physical sheet3.xml) is now named "Equity" and will hold the equity fund rows.

# ---------------------------------------------------------------------------
# 2) physical sheet1.xml ("Bond" tab) - 4 data rows
# ---------------------------------------------------------------------------
$wsBond = $wsEquitySlot

$wsBond.Range("A2:J8").ClearContents()

$wsBond.Range("A2").Value = "FLIA"
$wsBond.Range("B2").Value = "Traditional"
$wsBond.Range("C2").Value = "Franklin International Aggregate Bond ETF"
$wsBond.Range("D2").Value = "Cboe US"
$wsBond.Range("E2").Value = ""
$wsBond.Range("F2").Value = 0.0025
$wsBond.Range("G2").Value = 0.1517
$wsBond.Range("H2").Value = 4.26

$wsBond.Range("A3").Value = "BILS"
$wsBond.Range("B3").Value = "Traditional"
$wsBond.Range("C3").Value = "SPDR Bloomberg 3-12 Month T-Bill ETF"
$wsBond.Range("D3").Value = "NYSEArca"
$wsBond.Range("E3").Value = "Ultrashort Bond"
$wsBond.Range("F3").Value = 0.00135
$wsBond.Range("G3").Value = 0.0457
$wsBond.Range("H3").Value = 0

$wsBond.Range("A4").Value = "VCLT"
$wsBond.Range("B4").Value = "Traditional"
$wsBond.Range("C4").Value = "Vanguard Long-Term Corporate Bond Index Fund"
$wsBond.Range("D4").Value = "NasdaqGM"
$wsBond.Range("E4").Value = "Long-Term Bond"
$wsBond.Range("F4").Value = 0.0007
$wsBond.Range("G4").Value = 0.0535
$wsBond.Range("H4").Value = 13.98

$wsBond.Range("A5").Value = "VWOB"
$wsBond.Range("B5").Value = "Traditional"
$wsBond.Range("C5").Value = "Vanguard Emerging Markets Government Bond Index Fund"
$wsBond.Range("D5").Value = "NasdaqGM"
$wsBond.Range("E5").Value = "Emerging Markets Bond"
$wsBond.Range("F5").Value = 0.002
$wsBond.Range("G5").Value = 0.0597
$wsBond.Range("H5").Value = 11.02

# rows 6-8 no longer exist on the Bond tab (it only has 4 funds) - already
# cleared above via Range("A2:J8").ClearContents(), which also shrinks the
# sheet's used range/dimension down to A1:J5.

# ---------------------------------------------------------------------------
# 3) physical sheet2.xml ("Alternative" tab) - same rows/tickers, only the
#    Expense Ratio (F) and Dividend Yield (G) columns are rescaled from
#    percent-like numbers to decimal fractions.
# ---------------------------------------------------------------------------
$wsAlternative.Range("F2").Value = 0.0025
$wsAlternative.Range("G2").Value = 0

$wsAlternative.Range("F3").Value = 0.008500000000000001
$wsAlternative.Range("G3").Value = 0.0781

$wsAlternative.Range("F4").Value = 0.0067
$wsAlternative.Range("G4").Value = 0.025

$wsAlternative.Range("F5").Value = 0.009299999999999999
$wsAlternative.Range("G5").Value = 0

# ---------------------------------------------------------------------------
# 4) physical sheet3.xml ("Equity" tab) - 7 data rows
# ---------------------------------------------------------------------------
$wsEquity = $wsBondSlot

$wsEquity.Range("A2").Value = "VOO"
$wsEquity.Range("B2").Value = "Traditional"
$wsEquity.Range("C2").Value = "Vanguard 500 Index Fund"
$wsEquity.Range("D2").Value = "NYSEArca"
$wsEquity.Range("E2").Value = "Large Blend"
$wsEquity.Range("F2").Value = 0.0014
$wsEquity.Range("G2").Value = 0.0162
$wsEquity.Range("H2").Value = 18.7

$wsEquity.Range("A3").Value = "FLCA"
$wsEquity.Range("B3").Value = "Traditional"
$wsEquity.Range("C3").Value = "Franklin FTSE Canada ETF"
$wsEquity.Range("D3").Value = "NYSEArca"
$wsEquity.Range("E3").Value = "Miscellaneous Region"
$wsEquity.Range("F3").Value = 0.0009
$wsEquity.Range("G3").Value = 0.0317
$wsEquity.Range("H3").Value = 20.54

$wsEquity.Range("A4").Value = "FLJP"
$wsEquity.Range("B4").Value = "Traditional"
$wsEquity.Range("C4").Value = "Franklin FTSE Japan ETF"
$wsEquity.Range("D4").Value = "NYSEArca"
$wsEquity.Range("E4").Value = "Japan Stock"
$wsEquity.Range("F4").Value = 0.0009
$wsEquity.Range("G4").Value = 0.0224
$wsEquity.Range("H4").Value = 15.34

$wsEquity.Range("A5").Value = "FLAU"
$wsEquity.Range("B5").Value = "Traditional"
$wsEquity.Range("C5").Value = "Franklin FTSE Australia ETF"
$wsEquity.Range("D5").Value = "NYSEArca"
$wsEquity.Range("E5").Value = "Miscellaneous Region"
$wsEquity.Range("F5").Value = 0.0009
$wsEquity.Range("G5").Value = 0.0472
$wsEquity.Range("H5").Value = 23.72

$wsEquity.Range("A6").Value = "FLKR"
$wsEquity.Range("B6").Value = "Traditional"
$wsEquity.Range("C6").Value = "Franklin FTSE South Korea ETF"
$wsEquity.Range("D6").Value = "NYSEArca"
$wsEquity.Range("E6").Value = "Miscellaneous Region"
$wsEquity.Range("F6").Value = 0.0009
$wsEquity.Range("G6").Value = 0.0375
$wsEquity.Range("H6").Value = 25.41
$wsEquity.Range("I6").Value = ""
$wsEquity.Range("J6").Value = ""

$wsEquity.Range("A7").Value = "SPEU"
$wsEquity.Range("B7").Value = "Traditional"
$wsEquity.Range("C7").Value = "SPDR Portfolio Europe ETF"
$wsEquity.Range("D7").Value = "NYSEArca"
$wsEquity.Range("E7").Value = "Europe Stock"
$wsEquity.Range("F7").Value = 0.0007
$wsEquity.Range("G7").Value = 0.0333
$wsEquity.Range("H7").Value = 19.53
$wsEquity.Range("I7").Value = ""
$wsEquity.Range("J7").Value = ""

$wsEquity.Range("A8").Value = "SPEM"
$wsEquity.Range("B8").Value = "Traditional"
$wsEquity.Range("C8").Value = "SPDR Portfolio Emerging Markets ETF"
$wsEquity.Range("D8").Value = "NYSEArca"
$wsEquity.Range("E8").Value = "Diversified Emerging Mkts"
$wsEquity.Range("F8").Value = 0.0007
$wsEquity.Range("G8").Value = 0.0339
$wsEquity.Range("H8").Value = 17.73
$wsEquity.Range("I8").Value = ""
$wsEquity.Range("J8").Value = ""
